$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that need to be swapped between row 14 and row 15
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "S")

foreach ($col in $cols) {
    $addr14 = "$col" + "14"
    $addr15 = "$col" + "15"

    $val14 = $ws.Range($addr14).Value()
    $val15 = $ws.Range($addr15).Value()

    $ws.Range($addr14).Value = $val15
    $ws.Range($addr15).Value = $val14
}
